$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert a new record at row 31 (new market day), which
# pushes all the existing rows 31-60 down to 32-61 (preserving formatting,
# i.e. the date style on column D carries through the shift).
$ws.Rows("31").Insert()

# Populate the newly inserted row 31 with the new week's record.
$ws.Range("A31").Value = 4
$ws.Range("B31").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C31").Value = "Los Lagos"
$ws.Range("D31").Value = 44879
$ws.Range("E31").Value = 10
$ws.Range("F31").Value = 300000000
$ws.Range("G31").Value = "Espárragos"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 1500
$ws.Range("L31").Value = 1700
$ws.Range("M31").Value = 1600
$ws.Range("N31").Value = "$/kilo"
$ws.Range("O31").Value = "Provincia de Linares"
$ws.Range("P31").Value = 1600
$ws.Range("Q31").Value = 1
$ws.Range("R31").Value = "Hortaliza"
